$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8699934474986719
$ws.Range("C2").Value = 0.09993956321308417
$ws.Range("D2").Value = 0.01606691811218752
$ws.Range("E2").Value = 0.09030836886861948
$ws.Range("F2").Value = 5.465644298827414
$ws.Range("J2").Value = 0.2594876274752664
$ws.Range("K2").Value = 0.8912979301012172
$ws.Range("B3").Value = 0.8491418512583948
$ws.Range("C3").Value = 0.09828364450103777
$ws.Range("D3").Value = 0.01618851453896752
$ws.Range("E3").Value = 0.09047103471657181
$ws.Range("F3").Value = 5.267556758085533
$ws.Range("J3").Value = 0.2525470512218391
$ws.Range("K3").Value = 0.8712170596348301
$ws.Range("B4").Value = 0.8371533906245361
$ws.Range("C4").Value = 0.09740540495576511
$ws.Range("D4").Value = 0.01631623572598784
$ws.Range("E4").Value = 0.09065840970488281
$ws.Range("F4").Value = 5.146572118990235
$ws.Range("J4").Value = 0.2484037763694573
$ws.Range("K4").Value = 0.8597948119227965
$ws.Range("B5").Value = 0.8324722751297315
$ws.Range("C5").Value = 0.09708217566171129
$ws.Range("D5").Value = 0.01638145929810975
$ws.Range("E5").Value = 0.09075674096207464
$ws.Range("F5").Value = 5.097426765661965
$ws.Range("J5").Value = 0.246744946562103
$ws.Range("K5").Value = 0.8553676606779845
$ws.Range("B6").Value = 0.8317073021643466
$ws.Range("C6").Value = 0.09703059227260269
$ws.Range("D6").Value = 0.01639308057799482
$ws.Range("E6").Value = 0.09077439521360731
$ws.Range("F6").Value = 5.089275559796164
$ws.Range("J6").Value = 0.246471281852024
$ws.Range("K6").Value = 0.8546462567320816
$ws.Range("B7").Value = 0.8370894330473675
$ws.Range("C7").Value = 0.09740090565263415
$ws.Range("D7").Value = 0.0163170622239619
$ws.Range("E7").Value = 0.09065964689700934
$ws.Range("F7").Value = 5.145908698422119
$ws.Range("J7").Value = 0.2483812852033367
$ws.Range("K7").Value = 0.8597341855507352
$ws.Range("B8").Value = 0.8626344308604246
$ws.Range("C8").Value = 0.09933975414939766
$ws.Range("D8").Value = 0.01609773145044713
$ws.Range("E8").Value = 0.09034627862544653
$ws.Range("F8").Value = 5.39720768037958
$ws.Range("J8").Value = 0.2570698762938122
$ws.Range("K8").Value = 0.8841852675653001
$ws.Range("B9").Value = 0.9192199988744676
$ws.Range("C9").Value = 0.104248904238915
$ws.Range("D9").Value = 0.01609618782162414
$ws.Range("E9").Value = 0.09042750438570835
$ws.Range("F9").Value = 5.895328064028433
$ws.Range("J9").Value = 0.2750545855680997
$ws.Range("K9").Value = 0.9393705433123785
$ws.Range("B10").Value = 0.9647988458482075
$ws.Range("C10").Value = 0.1085426937193006
$ws.Range("D10").Value = 0.01636684744824635
$ws.Range("E10").Value = 0.09091373834475291
$ws.Range("F10").Value = 6.264898806538895
$ws.Range("J10").Value = 0.2888577489139266
$ws.Range("K10").Value = 0.9843847241739923
$ws.Range("B11").Value = 0.9864146321402245
$ws.Range("C11").Value = 0.1106479643102176
$ws.Range("D11").Value = 0.01655128841471054
$ws.Range("E11").Value = 0.0912281300920732
$ws.Range("F11").Value = 6.433887887222681
$ws.Range("J11").Value = 0.2952681591439159
$ws.Range("K11").Value = 1.005846633160274
$ws.Range("B12").Value = 0.9947275553084864
$ws.Range("C12").Value = 0.1114672482110137
$ws.Range("D12").Value = 0.01663013705012872
$ws.Range("E12").Value = 0.09136062820127933
$ws.Range("F12").Value = 6.498010456414306
$ws.Range("J12").Value = 0.2977147015141242
$ws.Range("K12").Value = 1.014116281974168
$ws.Range("B13").Value = 0.9929315368532343
$ws.Range("C13").Value = 0.1112898161571252
$ws.Range("D13").Value = 0.01661275208630997
$ws.Range("E13").Value = 0.09133149370794058
$ws.Range("F13").Value = 6.484194647867582
$ws.Range("J13").Value = 0.2971869442579589
$ws.Range("K13").Value = 1.012328913499431
$ws.Range("B14").Value = 0.9870959820739245
$ws.Range("C14").Value = 0.1107149239079632
$ws.Range("D14").Value = 0.01655759378832045
$ws.Range("E14").Value = 0.09123876108842666
$ws.Range("F14").Value = 6.439160653637373
$ws.Range("J14").Value = 0.2954690546328607
$ws.Range("K14").Value = 1.006524121277437
$ws.Range("B15").Value = 0.9835381609125022
$ws.Range("C15").Value = 0.1103656653197618
$ws.Range("D15").Value = 0.01652498603381503
$ws.Range("E15").Value = 0.09118371187832608
$ws.Range("F15").Value = 6.41159312177723
$ws.Range("J15").Value = 0.2944192851355893
$ws.Range("K15").Value = 1.002987101509603
$ws.Range("B16").Value = 0.9634040056862716
$ws.Range("C16").Value = 0.1084081855212133
$ws.Range("D16").Value = 0.01635604363829657
$ws.Range("E16").Value = 0.09089507115395179
$ws.Range("F16").Value = 6.253872898396452
$ws.Range("J16").Value = 0.2884414706379346
$ws.Range("K16").Value = 0.9830020284994703
$ws.Range("B17").Value = 0.9512786805744895
$ws.Range("C17").Value = 0.107246411866825
$ws.Range("D17").Value = 0.01626824049650111
$ws.Range("E17").Value = 0.09074190028689344
$ws.Range("F17").Value = 6.157342912812368
$ws.Range("J17").Value = 0.2848080270489675
$ws.Range("K17").Value = 0.9709946701608203
$ws.Range("B18").Value = 0.9443874646821087
$ws.Range("C18").Value = 0.1065924851834978
$ws.Range("D18").Value = 0.01622349853039395
$ws.Range("E18").Value = 0.09066257037889258
$ws.Range("F18").Value = 6.101902953871473
$ws.Range("J18").Value = 0.2827305122227983
$ws.Range("K18").Value = 0.9641809884463157
$ws.Range("B19").Value = 0.9420684431741222
$ws.Range("C19").Value = 0.1063735257958029
$ws.Range("D19").Value = 0.01620933328191398
$ws.Range("E19").Value = 0.09063721560814386
$ws.Range("F19").Value = 6.083145787751306
$ws.Range("J19").Value = 0.2820292163579978
$ws.Range("K19").Value = 0.9618898761130481
$ws.Range("B20").Value = 0.9525608520710875
$ws.Range("C20").Value = 0.107368603724467
$ws.Range("D20").Value = 0.01627698979125114
$ws.Range("E20").Value = 0.09075729763082663
$ws.Range("F20").Value = 6.167610219651607
$ws.Range("J20").Value = 0.2851935342774965
$ws.Range("K20").Value = 0.9722632804258637
$ws.Range("B21").Value = 0.9888065601286371
$ws.Range("C21").Value = 0.1108831832363961
$ws.Range("D21").Value = 0.0165735492242689
$ws.Range("E21").Value = 0.09126563367223994
$ws.Range("F21").Value = 6.452384658763719
$ws.Range("J21").Value = 0.2959731217371768
$ws.Range("K21").Value = 1.008225255150762
$ws.Range("B22").Value = 1.013238698459986
$ws.Range("C22").Value = 0.1133088664113018
$ws.Range("D22").Value = 0.01681994346732552
$ws.Range("E22").Value = 0.0916762480279445
$ws.Range("F22").Value = 6.639261763646573
$ws.Range("J22").Value = 0.3031293833486615
$ws.Range("K22").Value = 1.032559451972503
$ws.Range("B23").Value = 1.000130541618404
$ws.Range("C23").Value = 0.1120023882721028
$ws.Range("D23").Value = 0.01668356509190971
$ws.Range("E23").Value = 0.09144990848897194
$ws.Range("F23").Value = 6.539450662569266
$ws.Range("J23").Value = 0.299299718228923
$ws.Range("K23").Value = 1.019495500402599
$ws.Range("B24").Value = 0.9519809338463858
$ws.Range("C24").Value = 0.1073133172176455
$ws.Range("D24").Value = 0.01627301637889644
$ws.Range("E24").Value = 0.09075030930255679
$ws.Range("F24").Value = 6.162968194887128
$ws.Range("J24").Value = 0.2850192109637817
$ws.Range("K24").Value = 0.9716894629922308
$ws.Range("B25").Value = 0.9032122313031152
$ws.Range("C25").Value = 0.1028010900067216
$ws.Range("D25").Value = 0.01604964208171822
$ws.Range("E25").Value = 0.09033080284452311
$ws.Range("F25").Value = 5.759966764273514
$ws.Range("J25").Value = 0.270086611577355
$ws.Range("K25").Value = 0.92366082240315
